# The filtering threshold used to compute pre/recall/F1/PSC/PRC (and, for
# epsilon >= 0.066, NRC) was adjusted. Apply the updated metric values to
# the "EDCR Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-66 (epsilon 0.001 .. 0.065): NSC/NRC stay at 0, only pre/F1/PSC/PRC move.
$ws.Range("B2:B66").Value = 0.3023255813953488
$ws.Range("D2:D66").Value = 0.3537414965986394
$ws.Range("F2:F66").Value = 58
$ws.Range("H2:H66").Value = 1

# Rows 67-100 (epsilon 0.066 .. 0.099): NSC stays at 28, pre/F1/PSC/NRC/PRC move.
$ws.Range("B67:B100").Value = 0.3023255813953488
$ws.Range("D67:D100").Value = 0.3537414965986394
$ws.Range("F67:F100").Value = 86
$ws.Range("G67:G100").Value = 6
$ws.Range("H67:H100").Value = 1
